$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037192036402459
$ws.Range("D2").Value = 1.037920407686496
$ws.Range("E2").Value = 1.035928576258067
$ws.Range("F2").Value = 1.035890701401351
$ws.Range("I2").Value = 1.037128369315202
$ws.Range("J2").Value = 1.042296246904896
$ws.Range("K2").Value = 1.040709916208559
$ws.Range("L2").Value = 1.038723776063301
$ws.Range("M2").Value = 1.038686009657163
$ws.Range("N2").Value = 1.043776427128607

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038854277193227
$ws.Range("D3").Value = 1.038658757216885
$ws.Range("E3").Value = 1.037366102321112
$ws.Range("F3").Value = 1.038190592407985
$ws.Range("I3").Value = 1.037513795937972
$ws.Range("J3").Value = 1.043599183870308
$ws.Range("K3").Value = 1.04125836193355
$ws.Range("L3").Value = 1.039969131365413
$ws.Range("M3").Value = 1.040791436259562
$ws.Range("N3").Value = 1.045081214413903

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039926400807362
$ws.Range("D4").Value = 1.03913520772315
$ws.Range("E4").Value = 1.038293399558383
$ws.Range("F4").Value = 1.039674711336504
$ws.Range("I4").Value = 1.037760767237491
$ws.Range("J4").Value = 1.04443852953365
$ws.Range("K4").Value = 1.041611269572601
$ws.Range("L4").Value = 1.040771578285673
$ws.Range("M4").Value = 1.042149418448723
$ws.Range("N4").Value = 1.045921752044361

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040376310521044
$ws.Range("D5").Value = 1.039335197309797
$ws.Range("E5").Value = 1.038682559577106
$ws.Range("F5").Value = 1.04029769016468
$ws.Range("I5").Value = 1.037864017310145
$ws.Range("J5").Value = 1.044790508249617
$ws.Range("K5").Value = 1.041759163343974
$ws.Range("L5").Value = 1.041108129969889
$ws.Range("M5").Value = 1.042719294084363
$ws.Range("N5").Value = 1.046274230610478

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040451805261087
$ws.Range("D6").Value = 1.039368758355715
$ws.Range("E6").Value = 1.038747862015092
$ws.Range("F6").Value = 1.040402236464137
$ws.Range("I6").Value = 1.037881319751397
$ws.Range("J6").Value = 1.04484955568331
$ws.Range("K6").Value = 1.041783967990515
$ws.Range("L6").Value = 1.041164592101764
$ws.Range("M6").Value = 1.042814919643821
$ws.Range("N6").Value = 1.046333361898299

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039932415687212
$ws.Range("D7").Value = 1.03913788120755
$ws.Range("E7").Value = 1.038298602170615
$ws.Range("F7").Value = 1.039683039284258
$ws.Range("I7").Value = 1.03776214913038
$ws.Range("J7").Value = 1.044443236138441
$ws.Range("K7").Value = 1.041613247572345
$ws.Range("L7").Value = 1.040776078420403
$ws.Range("M7").Value = 1.042157037127907
$ws.Range("N7").Value = 1.045926465333071

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037754524389517
$ws.Range("D8").Value = 1.038170210037681
$ws.Range("E8").Value = 1.036414998942357
$ws.Range("F8").Value = 1.036668818972795
$ws.Range("I8").Value = 1.037259130689616
$ws.Range("J8").Value = 1.042737363500929
$ws.Range("K8").Value = 1.040895676946676
$ws.Range("L8").Value = 1.03914535737421
$ws.Range("M8").Value = 1.039398468933642
$ws.Range("N8").Value = 1.044218170160747

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033889570986664
$ws.Range("D9").Value = 1.036454842516318
$ws.Range("E9").Value = 1.03307320336267
$ws.Range("F9").Value = 1.031324973954534
$ws.Range("I9").Value = 1.036353989860937
$ws.Range("J9").Value = 1.039702127940913
$ws.Range("K9").Value = 1.039615937252812
$ws.Range("L9").Value = 1.036245358380554
$ws.Range("M9").Value = 1.034502878069634
$ws.Range("N9").Value = 1.04117862421809

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031293568948738
$ws.Range("D10").Value = 1.035304163780574
$ws.Range("E10").Value = 1.030829265469391
$ws.Range("F10").Value = 1.027738830328581
$ws.Range("I10").Value = 1.035737693793637
$ws.Range("J10").Value = 1.037658094114577
$ws.Range("K10").Value = 1.038752250339361
$ws.Range("L10").Value = 1.034293440768739
$ws.Range("M10").Value = 1.031214207100559
$ws.Range("N10").Value = 1.039131687629267

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030164640480066
$ws.Range("D11").Value = 1.034804168773224
$ws.Range("E11").Value = 1.029853612398562
$ws.Range("F11").Value = 1.026179995581962
$ws.Range("I11").Value = 1.03546772039527
$ws.Range("J11").Value = 1.036767938143274
$ws.Range("K11").Value = 1.038375711142508
$ws.Range("L11").Value = 1.033443652199766
$ws.Range("M11").Value = 1.029783889411774
$ws.Range("N11").Value = 1.038240267534382

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029744558046557
$ws.Range("D12").Value = 1.034618181617205
$ws.Range("E12").Value = 1.029490591909932
$ws.Range("F12").Value = 1.025600037076963
$ws.Range("I12").Value = 1.035366967447403
$ws.Range("J12").Value = 1.036436515117447
$ws.Range("K12").Value = 1.038235458586449
$ws.Range("L12").Value = 1.033127296956842
$ws.Range("M12").Value = 1.029251626922211
$ws.Range("N12").Value = 1.037908373849844

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029834701343153
$ws.Range("D13").Value = 1.034658088642514
$ws.Range("E13").Value = 1.029568489339052
$ws.Range("F13").Value = 1.025724483149822
$ws.Range("I13").Value = 1.035388600774868
$ws.Range("J13").Value = 1.036507642027415
$ws.Range("K13").Value = 1.038265560925488
$ws.Range("L13").Value = 1.033195188426914
$ws.Range("M13").Value = 1.029365843846595
$ws.Range("N13").Value = 1.037979601768182

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030129931717058
$ws.Range("D14").Value = 1.034788800486586
$ws.Range("E14").Value = 1.029823617745968
$ws.Range("F14").Value = 1.026132075366518
$ws.Range("I14").Value = 1.035459401800046
$ws.Range("J14").Value = 1.036740558611907
$ws.Range("K14").Value = 1.038364125789686
$ws.Range("L14").Value = 1.033417516668136
$ws.Range("M14").Value = 1.029739912583216
$ws.Range("N14").Value = 1.03821284912094

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030311733137592
$ws.Range("D15").Value = 1.034869300875769
$ws.Range("E15").Value = 1.029980728188717
$ws.Range("F15").Value = 1.026383081056238
$ws.Range("I15").Value = 1.035502961843094
$ws.Range("J15").Value = 1.036883962411224
$ws.Range("K15").Value = 1.038424803141728
$ws.Range("L15").Value = 1.033554406424073
$ws.Range("M15").Value = 1.029970258007056
$ws.Range("N15").Value = 1.0383564565701

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031368388435523
$ws.Range("D16").Value = 1.035337309660593
$ws.Range("E16").Value = 1.030893930389274
$ws.Range("F16").Value = 1.027842155454975
$ws.Range("I16").Value = 1.035755545015796
$ws.Range("J16").Value = 1.037717062457675
$ws.Range("K16").Value = 1.038777185709685
$ws.Range("L16").Value = 1.034349740297815
$ws.Range("M16").Value = 1.031308997081782
$ws.Range("N16").Value = 1.039190739714176

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032029889464308
$ws.Range("D17").Value = 1.035630409053713
$ws.Range("E17").Value = 1.031465671946371
$ws.Range("F17").Value = 1.028755760578278
$ws.Range("I17").Value = 1.03591314672099
$ws.Range("J17").Value = 1.038238273322384
$ws.Range("K17").Value = 1.038997537909056
$ws.Range("L17").Value = 1.034847391021126
$ws.Range("M17").Value = 1.03214704314494
$ws.Range("N17").Value = 1.039712690758097

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032415266002591
$ws.Range("D18").Value = 1.035801201076208
$ws.Range("E18").Value = 1.031798773043272
$ws.Range("F18").Value = 1.029288073401217
$ws.Range("I18").Value = 1.036004773105245
$ws.Range("J18").Value = 1.038541798209282
$ws.Range("K18").Value = 1.039125819239377
$ws.Range("L18").Value = 1.035137219942997
$ws.Range("M18").Value = 1.032635255426033
$ws.Range("N18").Value = 1.040016646685148

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032546591051259
$ws.Range("D19").Value = 1.035859408408075
$ws.Range("E19").Value = 1.031912286920932
$ws.Range("F19").Value = 1.029469481394462
$ws.Range("I19").Value = 1.03603596459801
$ws.Range("J19").Value = 1.038645209985174
$ws.Range("K19").Value = 1.039169518216394
$ws.Range("L19").Value = 1.03523596955594
$ws.Range("M19").Value = 1.032801621405975
$ws.Range("N19").Value = 1.040120205317623

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031958964994772
$ws.Range("D20").Value = 1.035598979663467
$ws.Range("E20").Value = 1.031404369541676
$ws.Range("F20").Value = 1.028657799388061
$ws.Range("I20").Value = 1.03589626861915
$ws.Range("J20").Value = 1.038182402988282
$ws.Range("K20").Value = 1.038973921723422
$ws.Range("L20").Value = 1.034794043671936
$ws.Range("M20").Value = 1.032057191626255
$ws.Range("N20").Value = 1.039656741081711

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030043014549492
$ws.Range("D21").Value = 1.034750316508995
$ws.Range("E21").Value = 1.029748506024403
$ws.Range("F21").Value = 1.026012075789567
$ws.Range("I21").Value = 1.035438565755892
$ws.Range("J21").Value = 1.036671992136821
$ws.Range("K21").Value = 1.038335111676099
$ws.Range("L21").Value = 1.033352066155473
$ws.Range("M21").Value = 1.029629785943817
$ws.Range("N21").Value = 1.038144185273594

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02883404156505
$ws.Range("D22").Value = 1.034215182532386
$ws.Range("E22").Value = 1.028703806312501
$ws.Range("F22").Value = 1.024343160501747
$ws.Range("I22").Value = 1.035148051702283
$ws.Range("J22").Value = 1.035717818612393
$ws.Range("K22").Value = 1.037931212521029
$ws.Range("L22").Value = 1.03244134550139
$ws.Range("M22").Value = 1.028097900289511
$ws.Range("N22").Value = 1.037188656713317

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029475359251782
$ws.Range("D23").Value = 1.034499015326755
$ws.Range("E23").Value = 1.029257967525048
$ws.Range("F23").Value = 1.02522841183929
$ws.Range("I23").Value = 1.03530231995998
$ws.Range("J23").Value = 1.036224078029021
$ws.Range("K23").Value = 1.038145542430186
$ws.Range("L23").Value = 1.032924528883705
$ws.Range("M23").Value = 1.028910531032759
$ws.Range("N23").Value = 1.037695635076385

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031991014148539
$ws.Range("D24").Value = 1.035613181770956
$ws.Range("E24").Value = 1.031432070652281
$ws.Range("F24").Value = 1.028702065616792
$ws.Range("I24").Value = 1.035903896033857
$ws.Range("J24").Value = 1.038207649896699
$ws.Range("K24").Value = 1.038984593622483
$ws.Range("L24").Value = 1.034818150409164
$ws.Range("M24").Value = 1.032097793530805
$ws.Range("N24").Value = 1.039682023843634

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034892089131911
$ws.Range("D25").Value = 1.036899539928347
$ws.Range("E25").Value = 1.033939907916896
$ws.Range("F25").Value = 1.032710509185649
$ws.Range("I25").Value = 1.036590239603518
$ws.Range("J25").Value = 1.040490365953922
$ws.Range("K25").Value = 1.039948617675215
$ws.Range("L25").Value = 1.036998294438557
$ws.Range("M25").Value = 1.035772779620571
$ws.Range("N25").Value = 1.041967981619489
